# Add a team "record" (Wins / Losses / Ties) to the roster sheet.
# This extends the used range from A1:AC59 to A1:AF59 by adding three
# new columns: AD (Wins), AE (Losses), AF (Ties).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) -----------------------------------------------
# Copy the formatting of the last existing header cell (AC1, which uses
# the bold/centered/bordered header style) onto the three new header
# cells so they match the rest of the header row.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (rows 2-59) ---------------------------------------------
# Every player row gets the same team record values.
$firstRow = 2
$lastRow = 59

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 72  # AD -> Wins
    $ws.Cells.Item($r, 31).Value = 90  # AE -> Losses
    $ws.Cells.Item($r, 32).Value = 0   # AF -> Ties
}
